$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 8507
$ws1.Range("F8").Value = 657
$ws1.Range("F9").Value = 123
$ws1.Range("F13").Value = 3680
$ws1.Range("F14").Value = 262
$ws1.Range("F24").Value = 435
$ws1.Range("F27").Value = 155
$ws1.Range("F28").Value = 345
$ws1.Range("F29").Value = 58
$ws1.Range("F35").Value = 44
$ws1.Range("F36").Value = 78
$ws1.Range("F39").Value = 156

# Sheet "全部类型" (sheetId 4) - mirrors the same rows shifted by +1
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 8507
$ws4.Range("F9").Value = 657
$ws4.Range("F10").Value = 123
$ws4.Range("F15").Value = 3680
$ws4.Range("F16").Value = 262
$ws4.Range("F29").Value = 435
$ws4.Range("F32").Value = 155
$ws4.Range("F34").Value = 345
$ws4.Range("F35").Value = 58
$ws4.Range("F41").Value = 44
$ws4.Range("F42").Value = 78
$ws4.Range("F45").Value = 156
